# Updated cryptos list on Thu Oct 31 23:36:48 UTC 2024 with GitHub Actions
#
# Applies the new "Price" (column D) and "Volume(1h)" (column E) values
# for each coin row, as produced by the upstream scraper run.
#
# Notes:
#  - Column D/E cells are plain text (inlineStr) in the source workbook,
#    several of the new "Price" values look like plain numbers
#    (e.g. "577.99", "0.600") which Excel would otherwise silently coerce
#    into numeric cells (and, worse, normalise away trailing zeros).
#    Prefixing the assignment with a literal apostrophe forces Excel to
#    store it as text, and resetting the cell .Style back to "Normal"
#    afterwards clears the quote-prefix styling bit so the cell keeps its
#    original (default) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, [string]$Address, [string]$Value)
    $cell = $Sheet.Range($Address)
    $cell.Value = "'" + $Value
    $cell.Style = "Normal"
}

# row -> @(new Price (or $null if unchanged), new Volume(1h))
$rowUpdates = @{
    2  = @("70.601.34",  "  -2.74%  ")
    3  = @("2.524.20",   "  -5.24%  ")
    4  = @($null,        "  -0.04%  ")
    5  = @("577.99",     "  -3.30%  ")
    6  = @("169.44",     "  -3.69%  ")
    7  = @($null,        "  +0.09%  ")
    8  = @($null,        "  -2.42%  ")
    9  = @("2.523.69",   "  -5.22%  ")
    10 = @("0.162",      "  -3.94%  ")
    11 = @($null,        "  -0.63%  ")
    12 = @($null,        "  -3.71%  ")
    13 = @("4.85",       "  -3.07%  ")
    14 = @("2.986.38",   "  -5.30%  ")
    15 = @("70.502.38",  "  -2.61%  ")
    16 = @($null,        "  -3.04%  ")
    17 = @("25.18",      "  -4.19%  ")
    18 = @("2.527.65",   "  -4.88%  ")
    19 = @($null,        "  -7.39%  ")
    20 = @("7.64",       "  -6.67%  ")
    21 = @("360.54",     "  -2.94%  ")
    22 = @($null,        "  -5.62%  ")
    23 = @("1.99",       "  -5.62%  ")
    25 = @("69.47",      "  -3.54%  ")
    26 = @($null,        "  -6.34%  ")
    27 = @("9.16",       "  -6.92%  ")
    28 = @("2.654.38",   "  -4.83%  ")
    29 = @("0.993",      "  -0.72%  ")
    30 = @($null,        "  -5.80%  ")
    31 = @($null,        "  -3.49%  ")
    32 = @("485.77",     "  -2.31%  ")
    33 = @($null,        "  -0.44%  ")
    34 = @($null,        "  -3.30%  ")
    36 = @("156.18",     "  -3.72%  ")
    37 = @($null,        "  +1.39%  ")
    38 = @($null,        "  -4.46%  ")
    39 = @("18.90",      "  -0.24%  ")
    40 = @($null,        "  +0.00%  ")
    41 = @("4.79",       "  -4.61%  ")
    42 = @($null,        "  -3.59%  ")
    43 = @($null,        "  -6.60%  ")
    44 = @("1.20",       "  -13.06%  ")
    45 = @($null,        "  -8.02%  ")
    46 = @($null,        "  -2.40%  ")
    47 = @("144.01",     "  -8.18%  ")
    48 = @("3.55",       "  -5.10%  ")
    49 = @($null,        "  -5.51%  ")
    50 = @($null,        "  -6.70%  ")
    51 = @("0.600",      "  -1.03%  ")
}

foreach ($row in $rowUpdates.Keys) {
    $price  = $rowUpdates[$row][0]
    $volume = $rowUpdates[$row][1]
    $priceAddr  = "D" + $row
    $volumeAddr = "E" + $row

    if ($null -ne $price) {
        Set-TextValue $ws $priceAddr $price
    }
    if ($null -ne $volume) {
        Set-TextValue $ws $volumeAddr $volume
    }
}
